$wb = $excel.ActiveWorkbook

# The workbook contains duplicated event data in the "展览" and "全部类型"
# sheets. Update the "想去人数" (F column) values for two events in both
# sheets to keep them in sync.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 54
    $ws.Range("F5").Value = 81
}
